# Summarizing Ourselves final data collection
# Appends a "TAGS:" heading (bold, matching the Role/Name labels already in
# the doc) followed by the list of subject tags harvested for this
# interview, to the end of the document body.

$d = $word.ActiveDocument

$tags = @(
    "1 Ellangyaraq -- Becoming Aware [Childhood]",
    "1 Neqsuryaraq -- Fishing",
    "1 Elluarrluni Yuuciquq -- One Will Live Well",
    "1 Tan'gaurluut Nasaurluut-llu Allakarluteng -- Boys and Girls Separate",
    "1 Nepiaq -- Sod-house",
    "1 Anglicarillerkaq, Tukercaryaraq -- Child Rearing",
    "1 Kass'at Tekiteqerraallratni -- Early Western Contact",
    "1 Nepiaq, Enepiaq, Enpiaq -- Sod-house",
    "1 Calirpagyaraq -- Hard Work",
    "1 Akusrarun -- Mischief, Misconduct"
)

# Build the raw WordprocessingML for the new paragraphs ourselves (rather
# than typing through the Range/Font API) so each run gets exactly the
# formatting it needs: bold run+mark for the "TAGS:" heading, and plain
# (no rPr at all) runs for every tag line, matching the style already used
# elsewhere in the document (e.g. "Role 1:" / "Interviewee").

function XmlEscape([string]$s) {
    $s = $s -replace '&', '&amp;'
    $s = $s -replace '<', '&lt;'
    $s = $s -replace '>', '&gt;'
    return $s
}

$xml = "<w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>TAGS:</w:t></w:r></w:p>"

for ($i = 0; $i -lt $tags.Count; $i++) {
    $text = XmlEscape $tags[$i]
    if ($i -eq $tags.Count - 1) {
        $xml += "<w:p><w:r><w:lastRenderedPageBreak/><w:t>$text</w:t></w:r></w:p>"
    } else {
        $xml += "<w:p><w:r><w:t>$text</w:t></w:r></w:p>"
    }
}

$end = $d.Content.End
$insertionPoint = $d.Range($end, $end)
$insertionPoint.InsertXML($xml)

Write-Output "ParaCount=$($d.Paragraphs.Count)"
